$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1442.138
$ws.Range("I17").Value = 1140.75
$ws.Range("J17").Value = 1654.8823
$ws.Range("K17").Value = 3422.25
$ws.Range("L17").Value = 4964.6469
$ws.Range("M17").Value = -3254.25
$ws.Range("N17").Value = -5300.6469
$ws.Range("H51").Value = 5518.75
$ws.Range("I51").Value = 1990
$ws.Range("J51").Value = 6695
$ws.Range("K51").Value = 1990
$ws.Range("L51").Value = 6695
$ws.Range("M51").Value = -1506
$ws.Range("N51").Value = -7663
$ws.Range("H132").Value = 21828052
$ws.Range("I132").Value = 23350198
$ws.Range("J132").Value = 10633.333
$ws.Range("K132").Value = 70050594
$ws.Range("L132").Value = 31899.999
$ws.Range("M132").Value = -70048064
$ws.Range("N132").Value = -36959.999
$ws.Range("H138").Value = 2625.23
$ws.Range("I138").Value = 1436.8966
$ws.Range("J138").Value = 3110.6057
$ws.Range("K138").Value = 4310.6898
$ws.Range("L138").Value = 9331.8171
$ws.Range("M138").Value = 829.3101999999999
$ws.Range("N138").Value = -19611.8171

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 9870.875
$ws.Range("J6").Value = 10499.25
$ws.Range("L6").Value = 10499.25
$ws.Range("N6").Value = -10845.25
$ws.Range("H32").Value = 6624.6
$ws.Range("I32").Value = 3640.3662
$ws.Range("J32").Value = 13930.827
$ws.Range("K32").Value = 3640.3662
$ws.Range("L32").Value = 13930.827
$ws.Range("M32").Value = -3353.3662
$ws.Range("N32").Value = -14504.827
$ws.Range("H61").Value = 1856.36
$ws.Range("I61").Value = 1654.5625
$ws.Range("J61").Value = 2215.111
$ws.Range("K61").Value = 1654.5625
$ws.Range("L61").Value = 2215.111
$ws.Range("M61").Value = -1442.5625
$ws.Range("N61").Value = -2639.111
$ws.Range("H74").Value = 1705.6875
$ws.Range("I74").Value = 1130.8846
$ws.Range("J74").Value = 4196.5
$ws.Range("K74").Value = 1130.8846
$ws.Range("L74").Value = 4196.5
$ws.Range("M74").Value = -256.8846000000001
$ws.Range("N74").Value = -5944.5
$ws.Range("H77").Value = 1705.6875
$ws.Range("I77").Value = 1130.8846
$ws.Range("J77").Value = 4196.5
$ws.Range("K77").Value = 5654.423000000001
$ws.Range("L77").Value = 20982.5
$ws.Range("M77").Value = -1286.423000000001
$ws.Range("N77").Value = -29718.5
$ws.Range("H132").Value = 1689.6492
$ws.Range("I132").Value = 841.9783
$ws.Range("J132").Value = 5234.4546
$ws.Range("K132").Value = 2525.9349
$ws.Range("L132").Value = 15703.3638
$ws.Range("M132").Value = 4.065099999999802
$ws.Range("N132").Value = -20763.3638
$ws.Range("H136").Value = 1856.36
$ws.Range("I136").Value = 1654.5625
$ws.Range("J136").Value = 2215.111
$ws.Range("K136").Value = 4963.6875
$ws.Range("L136").Value = 6645.333
$ws.Range("M136").Value = -2413.6875
$ws.Range("N136").Value = -11745.333
$ws.Range("H139").Value = 43340.19
$ws.Range("J139").Value = 43340.19
$ws.Range("L139").Value = 43340.19
$ws.Range("N139").Value = -53620.19

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 2001.8889
$ws.Range("I8").Value = 1534.6666
$ws.Range("J8").Value = 2936.3333
$ws.Range("K8").Value = 1534.6666
$ws.Range("L8").Value = 2936.3333
$ws.Range("M8").Value = -1394.6666
$ws.Range("N8").Value = -3216.3333
$ws.Range("H125").Value = 41592.5
$ws.Range("J125").Value = 41592.5
$ws.Range("L125").Value = 41592.5
$ws.Range("N125").Value = -51432.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2307.2778
$ws.Range("I105").Value = 1310.909
$ws.Range("J105").Value = 3873
$ws.Range("K105").Value = 1310.909
$ws.Range("L105").Value = 3873
$ws.Range("M105").Value = 436.0909999999999
$ws.Range("N105").Value = -7367
$ws.Range("H123").Value = 39780
$ws.Range("J123").Value = 39780
$ws.Range("L123").Value = 39780
$ws.Range("N123").Value = -49580
$ws.Range("H134").Value = 3334.3333
$ws.Range("I134").Value = 3332.2
$ws.Range("J134").Value = 3342.3333
$ws.Range("K134").Value = 9996.599999999999
$ws.Range("L134").Value = 10026.9999
$ws.Range("M134").Value = -7461.599999999999
$ws.Range("N134").Value = -15096.9999
$ws.Range("H138").Value = 49853.332
$ws.Range("J138").Value = 49853.332
$ws.Range("L138").Value = 49853.332
$ws.Range("N138").Value = -60133.332
$ws.Range("H140").Value = 119666.25
$ws.Range("J140").Value = 119666.25
$ws.Range("L140").Value = 119666.25
$ws.Range("N140").Value = -130026.25
$ws.Range("H141").Value = 18233.334
$ws.Range("J141").Value = 18233.334
$ws.Range("L141").Value = 18233.334
$ws.Range("N141").Value = -28593.334

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 684.8461
$ws.Range("I17").Value = 682.0909
$ws.Range("J17").Value = 700
$ws.Range("K17").Value = 2046.2727
$ws.Range("L17").Value = 2100
$ws.Range("M17").Value = -1877.2727
$ws.Range("N17").Value = -2438
$ws.Range("H80").Value = 5762.125
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 6156.7144
$ws.Range("K80").Value = 9000
$ws.Range("L80").Value = 18470.1432
$ws.Range("M80").Value = -8064
$ws.Range("N80").Value = -20342.1432
$ws.Range("H83").Value = 5762.125
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 6156.7144
$ws.Range("K83").Value = 27000
$ws.Range("L83").Value = 55410.4296
$ws.Range("M83").Value = -22320
$ws.Range("N83").Value = -64770.4296
$ws.Range("H136").Value = 2884
$ws.Range("J136").Value = 3720
$ws.Range("L136").Value = 11160
$ws.Range("N136").Value = -21360

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 753
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H15").Value = 37191.8
$ws.Range("J15").Value = 37191.8
$ws.Range("L15").Value = 37191.8
$ws.Range("N15").Value = -37767.8
$ws.Range("H81").Value = 37191.8
$ws.Range("J81").Value = 37191.8
$ws.Range("L81").Value = 37191.8
$ws.Range("N81").Value = -39187.8
$ws.Range("H84").Value = 37191.8
$ws.Range("J84").Value = 37191.8
$ws.Range("L84").Value = 111575.4
$ws.Range("N84").Value = -121559.4
$ws.Range("H122").Value = 2491.7083
$ws.Range("I122").Value = 1830.8823
$ws.Range("J122").Value = 4096.5713
$ws.Range("K122").Value = 5492.6469
$ws.Range("L122").Value = 12289.7139
$ws.Range("M122").Value = -3042.6469
$ws.Range("N122").Value = -17189.7139
$ws.Range("H126").Value = 3869.9053
$ws.Range("I126").Value = 2798.6863
$ws.Range("J126").Value = 5111.5454
$ws.Range("K126").Value = 8396.0589
$ws.Range("L126").Value = 15334.6362
$ws.Range("M126").Value = -5926.0589
$ws.Range("N126").Value = -20274.6362
$ws.Range("H127").Value = 30323
$ws.Range("J127").Value = 30323
$ws.Range("L127").Value = 30323
$ws.Range("N127").Value = -40243
$ws.Range("H132").Value = 2025.0646
$ws.Range("I132").Value = 908.8
$ws.Range("J132").Value = 2556.6191
$ws.Range("K132").Value = 2726.4
$ws.Range("L132").Value = 7669.8573
$ws.Range("M132").Value = -196.3999999999996
$ws.Range("N132").Value = -12729.8573
$ws.Range("H140").Value = 42578
$ws.Range("J140").Value = 42578
$ws.Range("L140").Value = 42578
$ws.Range("N140").Value = -52938

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5024.0786
$ws.Range("I40").Value = 4653.7856
$ws.Range("J40").Value = 6752.1113
$ws.Range("K40").Value = 4653.7856
$ws.Range("L40").Value = 6752.1113
$ws.Range("M40").Value = -4517.7856
$ws.Range("N40").Value = -7024.1113
$ws.Range("H61").Value = 2198
$ws.Range("I61").Value = 2329.3333
$ws.Range("J61").Value = 2066.6667
$ws.Range("K61").Value = 2329.3333
$ws.Range("L61").Value = 2066.6667
$ws.Range("M61").Value = -2127.3333
$ws.Range("N61").Value = -2470.6667
$ws.Range("H113").Value = 2198
$ws.Range("I113").Value = 2329.3333
$ws.Range("J113").Value = 2066.6667
$ws.Range("K113").Value = 2329.3333
$ws.Range("L113").Value = 2066.6667
$ws.Range("M113").Value = -159.3332999999998
$ws.Range("N113").Value = -6406.6667
$ws.Range("H122").Value = 5950.7
$ws.Range("I122").Value = 4760.4
$ws.Range("J122").Value = 7141
$ws.Range("K122").Value = 14281.2
$ws.Range("L122").Value = 21423
$ws.Range("M122").Value = -11831.2
$ws.Range("N122").Value = -26323
$ws.Range("H132").Value = 5678.436
$ws.Range("I132").Value = 2046
$ws.Range("J132").Value = 7105.4644
$ws.Range("K132").Value = 6138
$ws.Range("L132").Value = 21316.3932
$ws.Range("M132").Value = -3608
$ws.Range("N132").Value = -26376.3932
$ws.Range("H136").Value = 3330.0322
$ws.Range("I136").Value = 1654.2632
$ws.Range("J136").Value = 5983.3335
$ws.Range("K136").Value = 4962.7896
$ws.Range("L136").Value = 17950.0005
$ws.Range("M136").Value = -2412.7896
$ws.Range("N136").Value = -23050.0005
$ws.Range("H139").Value = 46488
$ws.Range("J139").Value = 48110
$ws.Range("L139").Value = 48110
$ws.Range("N139").Value = -58390
$ws.Range("H140").Value = 74983.336
$ws.Range("J140").Value = 74983.336
$ws.Range("L140").Value = 74983.336
$ws.Range("N140").Value = -85343.336
$ws.Range("H141").Value = 41970.418
$ws.Range("J141").Value = 41970.418
$ws.Range("L141").Value = 41970.418
$ws.Range("N141").Value = -52330.418

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 41400
$ws.Range("J80").Value = 41400
$ws.Range("L80").Value = 41400
$ws.Range("N80").Value = -43396
$ws.Range("H83").Value = 41400
$ws.Range("J83").Value = 41400
$ws.Range("L83").Value = 124200
$ws.Range("N83").Value = -134184
$ws.Range("H122").Value = 6198.222
$ws.Range("I122").Value = 3976.8
$ws.Range("J122").Value = 8975
$ws.Range("K122").Value = 11930.4
$ws.Range("L122").Value = 26925
$ws.Range("M122").Value = -9480.400000000001
$ws.Range("N122").Value = -31825
$ws.Range("H138").Value = 53616.668
$ws.Range("J138").Value = 53616.668
$ws.Range("L138").Value = 53616.668
$ws.Range("N138").Value = -63896.668
$ws.Range("H139").Value = 40776.8
$ws.Range("J139").Value = 41225.832
$ws.Range("L139").Value = 41225.832
$ws.Range("N139").Value = -51505.832
$ws.Range("H140").Value = 53547.57
$ws.Range("J140").Value = 53547.57
$ws.Range("L140").Value = 53547.57
$ws.Range("N140").Value = -63907.57
$ws.Range("H141").Value = 43517.918
$ws.Range("J141").Value = 43517.918
$ws.Range("L141").Value = 43517.918
$ws.Range("N141").Value = -53877.918
